$wb = $excel.ActiveWorkbook

# Sheet "OFF" - Week 16 update (row 2, "H")
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 476
$wsOff.Range("C2").Value = 334
$wsOff.Range("D2").Value = 96
$wsOff.Range("E2").Value = 38

# Sheet "DEF" - Week 16 update (row 2, "H")
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 413
$wsDef.Range("C2").Value = 305
$wsDef.Range("D2").Value = 98
$wsDef.Range("E2").Value = 53
$wsDef.Range("F2").Value = 6
